$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$notes = $s.NotesPage
$shp = $notes.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$lines = @(
  '# Name  ',
  'C.V. Raman  ',
  '  ',
  '# Textbook  ',
  'Knight, College Physics: A Strategic Approach, 2nd Edition, Chapter 18, section 1',
  'Knight, Physics for Scientists and Engineers: A Strategic Approach with Modern Physics, 3rd Edition, Chapter 22 ',
  'Knight, Physics for Scientists and Engineers: A Strategic Approach with Modern Physics, 3rd Edition, Chapter 24',
  '  ',
  '# Contributors  ',
  '  ',
  '# Description  ',
  'Sir Chandrasekhara Venkata Raman was the first Indian Physicist to win the',
  'Novel prize for his work in light scattering. He was born on November 7th,',
  '1988 in Tamil Nadu. His father was a professor of physics and mathematics and',
  'his mother came from a family of scholars. In 1921, when he was returning to',
  'India from London, on the voyage, he became fascinated with the blue color of',
  'Mediterranean sea. He had some physics apparatus on him - a prism, a miniature',
  'spectroscope and a diffraction grating - which he used to prove that the sea',
  'was scattering sunlight by water molecule, which refutes Lord Rayleigh''s',
  'explanation that the color of the sea is a reflection of the color of the sky.',
  'But Rayleigh argued that the blue color of the sky was caused by scattering of',
  'sunlight by the air molecules. From Raman''s discovery, he was obsessed with',
  'the idea of light scattering. Analyzing the effect of scattering light by',
  'water was not easy, so Raman made visual observations of how scattered lights',
  'of different color passed through liquid. This visual observation is now',
  'called the Raman scattering and is the result of the Raman effect - the change',
  'in the wavelength of light that occurs when a light beam is deflected by',
  'molecules. Raman effect was an important discovery because after its',
  'discovery, physicists used this technique to study the vibration and rotation',
  'of molecules. By the late 1930s, the Raman effect was a principal method of',
  'chemical analysis because Raman spectroscopy could be easily used to solid,',
  'liquid, gas and aqueous solutions.',
  '  ',
  '# Sources  ',
  'C.V. Raman The Raman Effect - Landmark. (n.d.). Retrieved from',
  'https://www.acs.org/content/acs/en/education/whatischemistry/landmarks/ramaneffect.html',
  '',
  'Home. (n.d.). Retrieved from https://www.famousscientists.org/c-v-raman/ Who',
  'is C.V. Raman? Everything You Need to Know. (n.d.). Retrieved from',
  'https://www.thefamouspeople.com/profiles/c-v-raman-5318.php  ',
  '  ',
  '# Photo  ',
  'https://www.thefamouspeople.com/profiles/images/c-v-raman-2.jpg'
)

$tr.Text = [string]::Join("`n", $lines)
